# Small BOM and documentation changes
# - Fix "McMaser" typo -> "McMaster" on the "Full Device" sheet (rows 20-23, column L)
# - Reduce quantity of Silicone Sealing Washer (row 23) from 3 to 2
# - Add a new BOM line (row 25) for a 1/4" Washer
# - Update the saved sheet view / selection state for "Full Device"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full Device")

# --- Fix the "McMaser" -> "McMaster" typo used as the vendor name ---
$ws.Range("L20").Value = "McMaster"
$ws.Range("L21").Value = "McMaster"
$ws.Range("L22").Value = "McMaster"
$ws.Range("L23").Value = "McMaster"

# --- Row 23: Silicone Sealing Washer quantity 3 -> 2 ---
$ws.Range("C23").Value = 2

# --- Row 25: new BOM entry for a 1/4" Washer ---
$ws.Range("B25").Value = '1/4" Washer'
$ws.Range("C25").Value = 2
$ws.Range("D25").Formula = "=3.3/100"
$ws.Range("F25").Formula = "=D25"
$ws.Range("H25").Formula = "=D25"
$ws.Range("J25").Formula = "=D25"
$ws.Range("L25").Value = "McMaster"
$ws.Range("M25").Value = "90108A413"
$ws.Range("N25").Value = "N/A"
$ws.Range("O25").Value = "https://www.mcmaster.com/#90108a413/=14xie31"
$ws.Range("P25").Value = "Goes around silicon washer to add thickness for rivet to clamp onto"

# --- Update the sheet's saved scroll position / active selection ---
$ws.Activate()
$ws.Range("M1").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("M26").Select()
